$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2, $null, $null, '29.097.60', '  +1.20%  '),
    @(3, $null, $null, '1.923.97', '  +2.06%  '),
    @(4, $null, $null, $null, '  +0.44%  '),
    @(5, $null, $null, '326.24', '  +1.25%  '),
    @(6, $null, $null, $null, '  +0.42%  '),
    @(7, $null, $null, '0.4609', '  +1.00%  '),
    @(8, $null, $null, '0.3826', '  +0.88%  '),
    @(9, $null, $null, '0.07764', '  +0.75%  '),
    @(10, $null, $null, '0.9784', $null),
    @(11, $null, $null, '22.73', '  +3.68%  '),
    @(12, $null, $null, '1.962.36', '  +3.98%  '),
    @(13, $null, $null, '5.695', '  +0.72%  '),
    @(14, $null, $null, '6.972', '  +0.52%  '),
    @(15, $null, $null, '0.07067', '  +1.13%  '),
    @(16, $null, $null, $null, '  +0.53%  '),
    @(17, $null, $null, '84.35', '  +1.53%  '),
    @(18, $null, $null, '0.000009514', '  +0.49%  '),
    @(19, $null, $null, '16.77', '  +1.41%  '),
    @(20, $null, $null, '1.006', '  +0.40%  '),
    @(21, $null, $null, '29.096.40', '  +1.34%  '),
    @(22, $null, $null, '5.342', '  +0.46%  '),
    @(23, $null, $null, '10.98', '  +1.15%  '),
    @(24, $null, $null, '2.071', '  -0.74%  '),
    @(25, $null, $null, '158.10', '  +1.87%  '),
    @(26, $null, $null, '19.12', '  +0.86%  '),
    @(27, $null, $null, '5.659', '  +1.33%  '),
    @(28, $null, $null, '118.14', '  +1.07%  '),
    @(29, $null, $null, $null, '  +2.08%  '),
    @(30, $null, $null, '0.09345', $null),
    @(31, $null, $null, '0.8544', '  +1.60%  '),
    @(32, $null, $null, '5.117', '  +1.15%  '),
    @(33, $null, $null, '1.244', '  +0.78%  '),
    @(34, $null, $null, '3.024', '  +1.37%  '),
    @(35, $null, $null, '1.161', '  +2.01%  '),
    @(36, $null, $null, $null, '  +0.54%  '),
    @(37, $null, $null, '3.167', '  +17.74%  '),
    @(38, $null, $null, $null, '  +0.41%  '),
    @(39, $null, $null, $null, '  +1.16%  '),
    @(40, $null, $null, '7.501', '  +1.21%  '),
    @(41, $null, $null, '0.5517', '  +0.87%  '),
    @(42, $null, $null, '0.1758', '  +0.74%  '),
    @(43, $null, $null, '9.328', '  +2.27%  '),
    @(44, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.188', '  +6.37%  '),
    @(45, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.000002714', '  -8.77%  '),
    @(46, $null, $null, '0.5188', '  +0.96%  '),
    @(47, $null, $null, '0.06927', '  +1.98%  '),
    @(48, $null, $null, '11.22', '  -0.35%  '),
    @(49, $null, $null, '110.39', '  -1.02%  '),
    @(50, $null, $null, '1.769', '  +0.01%  '),
    @(51, $null, $null, $null, '  +0.50%  '),
)

foreach ($item in $changes) {
    $r = $item[0]
    $bVal = $item[1]
    $cVal = $item[2]
    $dVal = $item[3]
    $eVal = $item[4]

    if ($bVal -ne $null) {
        $ws.Cells.Item($r, 2).Value = $bVal
    }
    if ($cVal -ne $null) {
        $ws.Cells.Item($r, 3).Value = $cVal
    }
    if ($dVal -ne $null) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.Value = "'" + $dVal
        $dCell.Style = "Normal"
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($r, 5).Value = $eVal
    }
}
